$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(13, 8).Value = 0
$ws.Cells.Item(13, 10).Value = 0
$ws.Cells.Item(13, 12).Value = 0
$ws.Cells.Item(13, 14).ClearContents()

$ws.Cells.Item(100, 8).Value = 23811310
$ws.Cells.Item(100, 9).Value = 1592
$ws.Cells.Item(100, 11).Value = 1592
$ws.Cells.Item(100, 13).Value = -1051

$ws.Cells.Item(129, 8).Value = 736.64703
$ws.Cells.Item(129, 9).Value = 664.9231
$ws.Cells.Item(129, 10).Value = 969.75
$ws.Cells.Item(129, 11).Value = 1994.7693
$ws.Cells.Item(129, 12).Value = 2909.25
$ws.Cells.Item(129, 13).Value = 3005.2307
$ws.Cells.Item(129, 14).Value = -12909.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(14, 8).Value = 800
$ws.Cells.Item(14, 9).Value = 800
$ws.Cells.Item(14, 11).Value = 800
$ws.Cells.Item(14, 13).Value = -625

$ws.Cells.Item(15, 8).Value = 12000
$ws.Cells.Item(15, 10).Value = 12000
$ws.Cells.Item(15, 12).Value = 12000
$ws.Cells.Item(15, 14).Value = -12700

$ws.Cells.Item(32, 8).Value = 28718.72
$ws.Cells.Item(32, 9).Value = 6891.0835
$ws.Cells.Item(32, 11).Value = 6891.0835
$ws.Cells.Item(32, 13).Value = -6604.0835

$ws.Cells.Item(80, 8).Value = 22335.75
$ws.Cells.Item(80, 10).Value = 22335.75
$ws.Cells.Item(80, 12).Value = 22335.75
$ws.Cells.Item(80, 14).Value = -24331.75

$ws.Cells.Item(82, 8).Value = 21250
$ws.Cells.Item(82, 9).Value = 20000
$ws.Cells.Item(82, 10).Value = 21666.666
$ws.Cells.Item(82, 11).Value = 20000
$ws.Cells.Item(82, 12).Value = 21666.666
$ws.Cells.Item(82, 13).Value = -19639
$ws.Cells.Item(82, 14).Value = -22388.666

$ws.Cells.Item(83, 8).Value = 22335.75
$ws.Cells.Item(83, 10).Value = 22335.75
$ws.Cells.Item(83, 12).Value = 67007.25
$ws.Cells.Item(83, 14).Value = -76991.25

$ws.Cells.Item(85, 8).Value = 21250
$ws.Cells.Item(85, 9).Value = 20000
$ws.Cells.Item(85, 10).Value = 21666.666
$ws.Cells.Item(85, 11).Value = 20000
$ws.Cells.Item(85, 12).Value = 21666.666
$ws.Cells.Item(85, 13).Value = -18752
$ws.Cells.Item(85, 14).Value = -24162.666

$ws.Cells.Item(122, 8).Value = 2000.4736
$ws.Cells.Item(122, 9).Value = 1916
$ws.Cells.Item(122, 10).Value = 2317.25
$ws.Cells.Item(122, 11).Value = 5748
$ws.Cells.Item(122, 12).Value = 6951.75
$ws.Cells.Item(122, 13).Value = -3298
$ws.Cells.Item(122, 14).Value = -11851.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(8, 8).Value = 2334.6667
$ws.Cells.Item(8, 9).Value = 2334.6667
$ws.Cells.Item(8, 10).Value = 0
$ws.Cells.Item(8, 11).Value = 2334.6667
$ws.Cells.Item(8, 12).Value = 0
$ws.Cells.Item(8, 13).Value = -2194.6667
$ws.Cells.Item(8, 14).ClearContents()

$ws.Cells.Item(11, 8).Value = 7371.4287
$ws.Cells.Item(11, 9).Value = 3266.6667
$ws.Cells.Item(11, 10).Value = 32000
$ws.Cells.Item(11, 11).Value = 3266.6667
$ws.Cells.Item(11, 12).Value = 32000
$ws.Cells.Item(11, 13).Value = -3126.6667
$ws.Cells.Item(11, 14).Value = -32280

$ws.Cells.Item(14, 8).Value = 1237.5
$ws.Cells.Item(14, 9).Value = 500
$ws.Cells.Item(14, 10).Value = 1483.3334
$ws.Cells.Item(14, 11).Value = 500
$ws.Cells.Item(14, 12).Value = 1483.3334
$ws.Cells.Item(14, 13).Value = -328
$ws.Cells.Item(14, 14).Value = -1827.3334

$ws.Cells.Item(82, 8).Value = 23935.564
$ws.Cells.Item(82, 10).Value = 29793.295
$ws.Cells.Item(82, 12).Value = 29793.295
$ws.Cells.Item(82, 14).Value = -30559.295

$ws.Cells.Item(85, 8).Value = 23935.564
$ws.Cells.Item(85, 10).Value = 29793.295
$ws.Cells.Item(85, 12).Value = 29793.295
$ws.Cells.Item(85, 14).Value = -32445.295

$ws.Cells.Item(122, 8).Value = 29750
$ws.Cells.Item(122, 10).Value = 29750
$ws.Cells.Item(122, 12).Value = 29750
$ws.Cells.Item(122, 14).Value = -39550

$ws.Cells.Item(125, 8).Value = 50370
$ws.Cells.Item(125, 10).Value = 50370
$ws.Cells.Item(125, 12).Value = 50370
$ws.Cells.Item(125, 14).Value = -60210

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(12, 8).Value = 6000000
$ws.Cells.Item(12, 9).Value = 0
$ws.Cells.Item(12, 10).Value = 6000000
$ws.Cells.Item(12, 11).Value = 0
$ws.Cells.Item(12, 12).Value = 6000000
$ws.Cells.Item(12, 13).ClearContents()
$ws.Cells.Item(12, 14).Value = -6000340

$ws.Cells.Item(41, 8).Value = 15903
$ws.Cells.Item(41, 10).Value = 19753.75
$ws.Cells.Item(41, 12).Value = 19753.75
$ws.Cells.Item(41, 14).Value = -20609.75

$ws.Cells.Item(50, 8).Value = 8912.286
$ws.Cells.Item(50, 10).Value = 8912.286
$ws.Cells.Item(50, 12).Value = 8912.286
$ws.Cells.Item(50, 14).Value = -10162.286

$ws.Cells.Item(109, 8).Value = 10950
$ws.Cells.Item(109, 10).Value = 10950
$ws.Cells.Item(109, 12).Value = 10950
$ws.Cells.Item(109, 14).Value = -13030

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(9, 8).Value = 14499.5
$ws.Cells.Item(9, 9).Value = 3000
$ws.Cells.Item(9, 10).Value = 15777.223
$ws.Cells.Item(9, 11).Value = 9000
$ws.Cells.Item(9, 12).Value = 47331.669
$ws.Cells.Item(9, 13).Value = -8776
$ws.Cells.Item(9, 14).Value = -47779.669

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(2, 8).Value = 48.692307
$ws.Cells.Item(2, 9).Value = 28.625
$ws.Cells.Item(2, 10).Value = 80.8
$ws.Cells.Item(2, 11).Value = 28.625
$ws.Cells.Item(2, 12).Value = 80.8
$ws.Cells.Item(2, 13).Value = 84.375
$ws.Cells.Item(2, 14).Value = -306.8

$ws.Cells.Item(14, 8).Value = 35334.668
$ws.Cells.Item(14, 9).Value = 100004
$ws.Cells.Item(14, 10).Value = 3000
$ws.Cells.Item(14, 11).Value = 100004
$ws.Cells.Item(14, 12).Value = 3000
$ws.Cells.Item(14, 13).Value = -99836
$ws.Cells.Item(14, 14).Value = -3336

$ws.Cells.Item(122, 8).Value = 1868.5625
$ws.Cells.Item(122, 9).Value = 2008.8182
$ws.Cells.Item(122, 10).Value = 1560
$ws.Cells.Item(122, 11).Value = 6026.4546
$ws.Cells.Item(122, 12).Value = 4680
$ws.Cells.Item(122, 13).Value = -3576.4546
$ws.Cells.Item(122, 14).Value = -9580

$ws.Cells.Item(123, 8).Value = 23574
$ws.Cells.Item(123, 10).Value = 23574
$ws.Cells.Item(123, 12).Value = 23574
$ws.Cells.Item(123, 14).Value = -28474

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 2506.0625
$ws.Cells.Item(7, 9).Value = 1789.7
$ws.Cells.Item(7, 10).Value = 3700
$ws.Cells.Item(7, 11).Value = 1789.7
$ws.Cells.Item(7, 12).Value = 3700
$ws.Cells.Item(7, 13).Value = -1677.7
$ws.Cells.Item(7, 14).Value = -3924

$ws.Cells.Item(17, 8).Value = 1996
$ws.Cells.Item(17, 9).Value = 2000
$ws.Cells.Item(17, 10).Value = 1980
$ws.Cells.Item(17, 11).Value = 2000
$ws.Cells.Item(17, 12).Value = 1980
$ws.Cells.Item(17, 13).Value = -1830
$ws.Cells.Item(17, 14).Value = -2320

$ws.Cells.Item(19, 8).Value = 103
$ws.Cells.Item(19, 9).Value = 103
$ws.Cells.Item(19, 11).Value = 103
$ws.Cells.Item(19, 13).Value = 67

$ws.Cells.Item(40, 8).Value = 1838
$ws.Cells.Item(40, 9).Value = 1611.75
$ws.Cells.Item(40, 10).Value = 2200
$ws.Cells.Item(40, 11).Value = 1611.75
$ws.Cells.Item(40, 12).Value = 2200
$ws.Cells.Item(40, 13).Value = -1475.75
$ws.Cells.Item(40, 14).Value = -2472

$ws.Cells.Item(61, 8).Value = 2250.3157
$ws.Cells.Item(61, 9).Value = 1596.6154
$ws.Cells.Item(61, 11).Value = 1596.6154
$ws.Cells.Item(61, 13).Value = -1394.6154

$ws.Cells.Item(113, 8).Value = 2250.3157
$ws.Cells.Item(113, 9).Value = 1596.6154
$ws.Cells.Item(113, 11).Value = 1596.6154
$ws.Cells.Item(113, 13).Value = 573.3846000000001

$ws.Cells.Item(118, 8).Value = 31000
$ws.Cells.Item(118, 10).Value = 31000
$ws.Cells.Item(118, 12).Value = 31000
$ws.Cells.Item(118, 14).Value = -34314

$ws.Cells.Item(122, 8).Value = 2772.652
$ws.Cells.Item(122, 9).Value = 2136.9092
$ws.Cells.Item(122, 10).Value = 3355.4167
$ws.Cells.Item(122, 11).Value = 6410.7276
$ws.Cells.Item(122, 12).Value = 10066.2501
$ws.Cells.Item(122, 13).Value = -3960.7276
$ws.Cells.Item(122, 14).Value = -14966.2501

$ws.Cells.Item(126, 8).Value = 2506.0625
$ws.Cells.Item(126, 9).Value = 1789.7
$ws.Cells.Item(126, 10).Value = 3700
$ws.Cells.Item(126, 11).Value = 5369.1
$ws.Cells.Item(126, 12).Value = 11100
$ws.Cells.Item(126, 13).Value = -2899.1
$ws.Cells.Item(126, 14).Value = -16040

$ws.Cells.Item(127, 8).Value = 47038.09
$ws.Cells.Item(127, 10).Value = 47038.09
$ws.Cells.Item(127, 12).Value = 47038.09
$ws.Cells.Item(127, 14).Value = -56958.09

$ws.Cells.Item(132, 8).Value = 4260.7
$ws.Cells.Item(132, 9).Value = 4162.769
$ws.Cells.Item(132, 10).Value = 4442.5713
$ws.Cells.Item(132, 11).Value = 12488.307
$ws.Cells.Item(132, 12).Value = 13327.7139
$ws.Cells.Item(132, 13).Value = -9958.307000000001
$ws.Cells.Item(132, 14).Value = -18387.7139

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(12, 8).Value = 70007
$ws.Cells.Item(12, 10).Value = 70007
$ws.Cells.Item(12, 12).Value = 70007
$ws.Cells.Item(12, 14).Value = -70291

$ws.Cells.Item(75, 8).Value = 28032.5
$ws.Cells.Item(75, 10).Value = 28032.5
$ws.Cells.Item(75, 12).Value = 28032.5
$ws.Cells.Item(75, 14).Value = -29904.5

$ws.Cells.Item(78, 8).Value = 28032.5
$ws.Cells.Item(78, 10).Value = 28032.5
$ws.Cells.Item(78, 12).Value = 84097.5
$ws.Cells.Item(78, 14).Value = -93457.5

$ws.Cells.Item(100, 8).Value = 556573.9
$ws.Cells.Item(100, 9).Value = 1095.4375
$ws.Cells.Item(100, 11).Value = 2190.875
$ws.Cells.Item(100, 13).Value = -1649.875

$ws.Cells.Item(109, 8).Value = 25800
$ws.Cells.Item(109, 10).Value = 25800
$ws.Cells.Item(109, 12).Value = 25800
$ws.Cells.Item(109, 14).Value = -28574

$ws.Cells.Item(121, 8).Value = 29800
$ws.Cells.Item(121, 10).Value = 29800
$ws.Cells.Item(121, 12).Value = 29800
$ws.Cells.Item(121, 14).Value = -33294
